$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of graphics card data
$ws.Range("A40").Value = "ASUS Prime Radeon RX 9060 XT 16GB GDDR6 OC Edition"
$ws.Range("B40").Value = 107001
$ws.Range("C40").Value = 4711387994214

$ws.Range("A41").Value = "ASUS TUF Gaming Nvidia GeForce RTX 4070 Super"
$ws.Range("B41").Value = 106424
$ws.Range("C41").Value = 4711387450871

# Copy style from existing row (row 39) for A and B columns
$ws.Range("A39:B39").Copy()
$ws.Range("A40:B41").PasteSpecial(-4122)

# Update selection
$ws.Range("J28").Select()
